$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '72.541.75'
$ws.Cells.Item(2, 5).Value = '  +5.71%  '
Set-TextValue 3 4 '4.062.33'
$ws.Cells.Item(3, 5).Value = '  +5.94%  '
$ws.Cells.Item(4, 5).Value = '  -0.15%  '
Set-TextValue 5 4 '521.06'
$ws.Cells.Item(5, 5).Value = '  -0.03%  '
Set-TextValue 6 4 '147.68'
$ws.Cells.Item(6, 5).Value = '  +4.48%  '
$ws.Cells.Item(7, 5).Value = '  +21.06%  '
Set-TextValue 8 4 '4.054.98'
$ws.Cells.Item(8, 5).Value = '  +5.92%  '
$ws.Cells.Item(9, 5).Value = '  +0.13%  '
$ws.Cells.Item(10, 5).Value = '  +10.86%  '
$ws.Cells.Item(11, 5).Value = '  +5.85%  '
Set-TextValue 12 4 '0.0000332'
$ws.Cells.Item(12, 5).Value = '  +1.91%  '
Set-TextValue 13 4 '48.78'
$ws.Cells.Item(13, 5).Value = '  +17.86%  '
$ws.Cells.Item(14, 5).Value = '  +10.59%  '
Set-TextValue 15 4 '4.707.58'
$ws.Cells.Item(15, 5).Value = '  +5.87%  '
Set-TextValue 16 4 '4.075.03'
$ws.Cells.Item(16, 5).Value = '  +6.53%  '
$ws.Cells.Item(17, 2).Value = 'Uniswap'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 17 4 '14.50'
$ws.Cells.Item(17, 5).Value = '  +5.36%  '
$ws.Cells.Item(18, 2).Value = 'Chainlink'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 18 4 '21.42'
$ws.Cells.Item(18, 5).Value = '  +2.99%  '
$ws.Cells.Item(19, 5).Value = '  +3.07%  '
$ws.Cells.Item(20, 5).Value = '  -0.06%  '
Set-TextValue 21 4 '72.456.93'
$ws.Cells.Item(21, 5).Value = '  +5.54%  '
Set-TextValue 22 4 '447.99'
$ws.Cells.Item(22, 5).Value = '  +7.27%  '
Set-TextValue 23 4 '104.83'
$ws.Cells.Item(23, 5).Value = '  +21.37%  '
Set-TextValue 24 4 '3.60'
$ws.Cells.Item(24, 5).Value = '  +7.06%  '
Set-TextValue 25 4 '15.11'
$ws.Cells.Item(25, 5).Value = '  +8.26%  '
Set-TextValue 26 4 '4.02'
$ws.Cells.Item(26, 5).Value = '  +2.41%  '
Set-TextValue 27 4 '11.39'
$ws.Cells.Item(27, 5).Value = '  +1.31%  '
Set-TextValue 28 4 '11.12'
$ws.Cells.Item(28, 5).Value = '  +6.86%  '
Set-TextValue 29 4 '38.15'
$ws.Cells.Item(29, 5).Value = '  +6.78%  '
Set-TextValue 30 4 '5.84'
$ws.Cells.Item(30, 5).Value = '  +3.10%  '
Set-TextValue 31 4 '3.30'
$ws.Cells.Item(31, 5).Value = '  +16.82%  '
$ws.Cells.Item(32, 5).Value = '  +5.65%  '
Set-TextValue 33 4 '0.131'
$ws.Cells.Item(33, 5).Value = '  +5.18%  '
Set-TextValue 34 4 '679.51'
$ws.Cells.Item(34, 5).Value = '  +0.03%  '
Set-TextValue 35 4 '68.10'
$ws.Cells.Item(35, 5).Value = '  +1.57%  '
Set-TextValue 36 4 '6.62'
$ws.Cells.Item(36, 5).Value = '  +14.05%  '
Set-TextValue 37 4 '42.22'
$ws.Cells.Item(37, 5).Value = '  +7.47%  '
$ws.Cells.Item(38, 5).Value = '  +3.47%  '
Set-TextValue 39 4 '0.431'
$ws.Cells.Item(39, 5).Value = '  -0.83%  '
$ws.Cells.Item(40, 5).Value = '  +4.73%  '
$ws.Cells.Item(41, 5).Value = '  +9.93%  '
$ws.Cells.Item(42, 5).Value = '  +0.05%  '
Set-TextValue 43 4 '0.0500'
$ws.Cells.Item(43, 5).Value = '  +5.56%  '
Set-TextValue 44 4 '1.00'
$ws.Cells.Item(44, 5).Value = '  +0.03%  '
$ws.Cells.Item(45, 2).Value = 'WEMIXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 45 4 '3.21'
$ws.Cells.Item(45, 5).Value = '  +2.85%  '
$ws.Cells.Item(46, 2).Value = 'Stellar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 46 4 '0.159'
$ws.Cells.Item(46, 5).Value = '  +15.06%  '
$ws.Cells.Item(47, 2).Value = 'Fetch.AI'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 47 4 '2.69'
$ws.Cells.Item(47, 5).Value = '  -0.48%  '
$ws.Cells.Item(48, 2).Value = 'THORChain'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 48 4 '9.85'
$ws.Cells.Item(48, 5).Value = '  +17.37%  '
Set-TextValue 49 4 '3.46'
$ws.Cells.Item(49, 5).Value = '  +2.15%  '
$ws.Cells.Item(50, 5).Value = '  +5.80%  '
$ws.Cells.Item(51, 2).Value = 'FLOKI'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue 51 4 '0.000275'
$ws.Cells.Item(51, 5).Value = '  +3.15%  '
